# khl_stats_1369_ext.xlsx — publish files + archive (2025-11-16 11:06:41)
#
# 1) Matches_SOG: append 3 new match rows (474-476) for 2025-11-15 games.
# 2) Shots_HA / Shots_Summary: refresh as_of_utc timestamps for every team
#    row, and bump the shot counters for the teams that played on 2025-11-15.
# 3) Meta_ext: refresh as_of_utc + bump build_version.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Matches_SOG — new rows
# ---------------------------------------------------------------------
$matches = $wb.Worksheets.Item("Matches_SOG")

$newMatches = @(
    @{ Row = 474; Uid = "897773"; Date = "2025-11-15T13:30:00"; Home = "Сибирь";    Away = "Динамо М"; SogHome = 30; SogAway = 24 },
    @{ Row = 475; Uid = "897771"; Date = "2025-11-15T15:00:00"; Home = "ХК Сочи";   Away = "Драконы";  SogHome = 18; SogAway = 38 },
    @{ Row = 476; Uid = "897772"; Date = "2025-11-15T17:00:00"; Home = "Локомотив"; Away = "Барыс";    SogHome = 27; SogAway = 31 }
)

foreach ($m in $newMatches) {
    $r = $m.Row
    $matches.Cells.Item($r, 1).NumberFormat = "@"
    $matches.Cells.Item($r, 1).Value = $m.Uid
    $matches.Cells.Item($r, 2).NumberFormat = "@"
    $matches.Cells.Item($r, 2).Value = $m.Date
    $matches.Cells.Item($r, 3).Value = $m.Home
    $matches.Cells.Item($r, 4).Value = $m.Away
    $matches.Cells.Item($r, 5).Value = $m.SogHome
    $matches.Cells.Item($r, 6).Value = $m.SogAway
    $matches.Cells.Item($r, 7).Value = "khl_text"
}

# ---------------------------------------------------------------------
# 2) Shots_HA — as_of_utc refresh (D2:D23) + counter bumps
# ---------------------------------------------------------------------
$shotsHA = $wb.Worksheets.Item("Shots_HA")

for ($r = 2; $r -le 23; $r++) {
    $shotsHA.Cells.Item($r, 4).Value = "2025-11-15T17:00:00Z"
}

# Барыс (row 7)
$shotsHA.Range("F7").Value = 16
$shotsHA.Range("K7").Value = 446
$shotsHA.Range("L7").Value = 542
$shotsHA.Range("M7").Value = 27.9
$shotsHA.Range("N7").Value = 33.9

# Динамо М (row 8)
$shotsHA.Range("F8").Value = 23
$shotsHA.Range("K8").Value = 647
$shotsHA.Range("L8").Value = 736
$shotsHA.Range("M8").Value = 28.1
$shotsHA.Range("N8").Value = 32

# Драконы (row 10)
$shotsHA.Range("F10").Value = 22
$shotsHA.Range("K10").Value = 611
$shotsHA.Range("L10").Value = 790
$shotsHA.Range("M10").Value = 27.8
$shotsHA.Range("N10").Value = 35.9

# Локомотив (row 12)
$shotsHA.Range("E12").Value = 19
$shotsHA.Range("G12").Value = 582
$shotsHA.Range("H12").Value = 526
$shotsHA.Range("I12").Value = 30.6
$shotsHA.Range("J12").Value = 27.7

# Сибирь (row 18)
$shotsHA.Range("E18").Value = 20
$shotsHA.Range("G18").Value = 541
$shotsHA.Range("H18").Value = 757
$shotsHA.Range("I18").Value = 27.1
$shotsHA.Range("J18").Value = 37.9

# ХК Сочи (row 22)
$shotsHA.Range("E22").Value = 20
$shotsHA.Range("G22").Value = 578
$shotsHA.Range("H22").Value = 660
$shotsHA.Range("I22").Value = 28.9
$shotsHA.Range("J22").Value = 33

# ---------------------------------------------------------------------
# 3) Shots_Summary — as_of_utc refresh (D2:D23) + counter bumps
# ---------------------------------------------------------------------
$shotsSummary = $wb.Worksheets.Item("Shots_Summary")

for ($r = 2; $r -le 23; $r++) {
    $shotsSummary.Cells.Item($r, 4).Value = "2025-11-15T17:00:00Z"
}

# Барыс (row 7)
$shotsSummary.Range("E7").Value = 46
$shotsSummary.Range("F7").Value = 1403
$shotsSummary.Range("G7").Value = 1495
$shotsSummary.Range("I7").Value = 32.5

# Динамо М (row 8)
$shotsSummary.Range("E8").Value = 41
$shotsSummary.Range("F8").Value = 1245
$shotsSummary.Range("G8").Value = 1227
$shotsSummary.Range("H8").Value = 30.4

# Драконы (row 10)
$shotsSummary.Range("E10").Value = 42
$shotsSummary.Range("F10").Value = 1182
$shotsSummary.Range("G10").Value = 1488
$shotsSummary.Range("H10").Value = 28.1
$shotsSummary.Range("I10").Value = 35.4

# Локомотив (row 12)
$shotsSummary.Range("E12").Value = 46
$shotsSummary.Range("F12").Value = 1421
$shotsSummary.Range("G12").Value = 1192
$shotsSummary.Range("H12").Value = 30.9
$shotsSummary.Range("I12").Value = 25.9

# Сибирь (row 18)
$shotsSummary.Range("E18").Value = 44
$shotsSummary.Range("F18").Value = 1198
$shotsSummary.Range("G18").Value = 1539
$shotsSummary.Range("I18").Value = 35

# ХК Сочи (row 22)
$shotsSummary.Range("E22").Value = 40
$shotsSummary.Range("F22").Value = 1094
$shotsSummary.Range("G22").Value = 1405
$shotsSummary.Range("H22").Value = 27.4

# ---------------------------------------------------------------------
# 4) Meta_ext — as_of_utc refresh + build_version bump
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Meta_ext")
$meta.Range("B2").Value = "2025-11-15T17:00:00Z"
$meta.Range("D2").Value = 67
